$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E (particip) and F (taxa_sucesso) hold values that were stored as
# fractions (e.g. 0.97) but should be stored as whole-number percentages
# (e.g. 97) while keeping the existing "0.00%" cell formatting.
for ($row = 2; $row -le 7; $row++) {
    foreach ($col in @("E", "F")) {
        $cell = $ws.Range("$col$row")
        $cell.Value2 = $cell.Value2 * 100
    }
}
